$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header text updates ----
$ws.Range("A8").Value = "Volume 32   Number  29"
$ws.Range("C9").Value = "Report Covering the Week  7/14/2025  Through  7/20/2025"

# ---- Column H width (bestFit autosize side-effect) ----
$ws.Columns.Item(8).ColumnWidth = $ws.Columns.Item(5).ColumnWidth

# ---- Row 16 ----
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 2
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = -30
$ws.Range("I16").Value = 60
$ws.Range("J16").Value = 88
$ws.Range("K16").Value = -31.818181818181
$ws.Range("L16").Value = -40
$ws.Range("M16").Value = -17.808219178082
$ws.Range("N16").Value = -86.870897155361

# ---- Row 17 ----
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 3
$ws.Range("F17").Value = 30
$ws.Range("G17").Value = 27
$ws.Range("H17").Value = 11.111111111111
$ws.Range("I17").Value = 77
$ws.Range("J17").Value = 94
$ws.Range("K17").Value = -18.085106382978
$ws.Range("L17").Value = -34.745762711864
$ws.Range("M17").Value = 57.142857142857
$ws.Range("N17").Value = -50.955414012738

# ---- Row 18 ----
$ws.Range("C18").Value = 6
$ws.Range("E18").Value = 20
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = 7.142857142857
$ws.Range("I18").Value = 90
$ws.Range("J18").Value = 145
$ws.Range("K18").Value = -37.931034482758
$ws.Range("L18").Value = -47.368421052631
$ws.Range("M18").Value = -15.094339622641
$ws.Range("N18").Value = -79.729729729729

# ---- Row 19 ----
$ws.Range("C19").Value = 18
$ws.Range("D19").Value = 26
$ws.Range("E19").Value = -30.76923076923
$ws.Range("F19").Value = 90
$ws.Range("G19").Value = 118
$ws.Range("H19").Value = -23.728813559322
$ws.Range("I19").Value = 516
$ws.Range("J19").Value = 595
$ws.Range("K19").Value = -13.277310924369
$ws.Range("L19").Value = -23.781388478582
$ws.Range("M19").Value = -10.104529616724
$ws.Range("N19").Value = -60.640732265446

# ---- Row 20 ----
$ws.Range("C20").Value = 1
$ws.Range("F20").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("D20").Value = 2
$ws.Range("F20").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("E20").Value = -50
$ws.Range("H20").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 3
$ws.Range("I20").Value = 9
$ws.Range("J20").Value = 23
$ws.Range("K20").Value = -60.869565217391
$ws.Range("L20").Value = -59.090909090909
$ws.Range("M20").Value = -59.090909090909
$ws.Range("N20").Value = -97.744360902255

# ---- Row 21 ----
$ws.Range("C21").Value = 30
$ws.Range("D21").Value = 38
$ws.Range("E21").Value = -21.052631578947
$ws.Range("F21").Value = 145
$ws.Range("G21").Value = 172
$ws.Range("H21").Value = -15.697674418604
$ws.Range("I21").Value = 758
$ws.Range("J21").Value = 947
$ws.Range("K21").Value = -19.957761351636
$ws.Range("L21").Value = -30.649588289112
$ws.Range("M21").Value = -8.67469879518
$ws.Range("N21").Value = -72.684684684684

# ---- Row 22 ----
$ws.Range("D22").Value = 1
$ws.Range("F22").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = -100
$ws.Range("H22").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 26
$ws.Range("K22").Value = 15.384615384615
$ws.Range("L22").Value = 20
$ws.Range("M22").Value = -9.090909090909

# ---- Row 24 ----
$ws.Range("C24").Value = 21
$ws.Range("D24").Value = 61
$ws.Range("E24").Value = -65.573770491803
$ws.Range("F24").Value = 115
$ws.Range("G24").Value = 160
$ws.Range("H24").Value = -28.125
$ws.Range("I24").Value = 791
$ws.Range("J24").Value = 952
$ws.Range("K24").Value = -16.911764705882
$ws.Range("L24").Value = -31.217391304347
$ws.Range("M24").Value = -6.941176470588

# ---- Row 25 ----
$ws.Range("C25").Value = 15
$ws.Range("D25").Value = 43
$ws.Range("E25").Value = -65.116279069767
$ws.Range("F25").Value = 83
$ws.Range("G25").Value = 117
$ws.Range("H25").Value = -29.059829059829
$ws.Range("I25").Value = 577
$ws.Range("J25").Value = 771
$ws.Range("K25").Value = -25.162127107652
$ws.Range("L25").Value = -33.601841196777

# ---- Row 26 ----
$ws.Range("C26").Value = 5
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = -28.571428571428
$ws.Range("F26").Value = 43
$ws.Range("H26").Value = 16.216216216216
$ws.Range("I26").Value = 211
$ws.Range("J26").Value = 207
$ws.Range("K26").Value = 1.932367149758
$ws.Range("L26").Value = -15.936254980079
$ws.Range("M26").Value = 57.462686567164

# ---- Row 27 ----
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range("C27").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "***.*"
$ws.Range("C27").Copy()
$ws.Range("E27").PasteSpecial(-4122)

# ---- Row 28 ----
$ws.Range("C28").Value = 5
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = 150
$ws.Range("F28").Value = 9
$ws.Range("G28").Value = 7
$ws.Range("H28").Value = 28.571428571428
$ws.Range("I28").Value = 45
$ws.Range("J28").Value = 39
$ws.Range("K28").Value = 15.384615384615
$ws.Range("L28").Value = 18.421052631578

# ---- Row 31 ----
$ws.Range("F31").NumberFormat = "@"
$ws.Range("F31").Value = "0"
$ws.Range("C31").Copy()
$ws.Range("F31").PasteSpecial(-4122)
$ws.Range("G31").Value = 1
$ws.Range("H31").Value = -100
$ws.Range("I31").Value = 7
$ws.Range("K31").Value = -58.823529411764
$ws.Range("L31").Value = 0

$excel.CutCopyMode = $false
Write-Host "Edit complete"